# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E on this sheet are stored as plain text (not numbers/percent-formatted
# numbers), so pre-format any price cell that would otherwise look like a plain
# number ("58.01", "1.00", ...) as Text before writing it, keeping it a literal string.
$textCells = @("D5", "D7", "D11", "D13", "D15", "D19", "D20", "D22", "D23", "D24", "D27", "D28", "D29", "D32", "D33", "D34", "D35", "D39", "D41", "D44", "D46", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.085.00"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.066.30"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "253.28"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").Value = "58.01"
$ws.Range("E7").Value = "  +5.98%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").Value = "0.0801"
$ws.Range("E11").Value = "  +7.15%  "
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("D13").Value = "16.34"
$ws.Range("E13").Value = "  +8.88%  "
$ws.Range("D14").Value = "2.368.87"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "0.809"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("E16").Value = "  +8.93%  "
$ws.Range("D17").Value = "2.068.34"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "37.069.99"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "16.57"
$ws.Range("E19").Value = "  +14.33%  "
$ws.Range("D20").Value = "75.65"
$ws.Range("E20").Value = "  +3.91%  "
$ws.Range("D21").Value = "0.0₃0925"
$ws.Range("E21").Value = "  +9.22%  "
$ws.Range("D22").Value = "5.47"
$ws.Range("E22").Value = "  +4.95%  "
$ws.Range("D23").Value = "238.83"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("E26").Value = "  +12.69%  "
$ws.Range("D27").Value = "169.51"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").Value = "9.35"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "20.35"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("D32").Value = "4.78"
$ws.Range("E32").Value = "  +5.63%  "
$ws.Range("D33").Value = "0.0621"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").Value = "4.50"
$ws.Range("E34").Value = "  +7.78%  "
$ws.Range("D35").Value = "0.0912"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").Value = "0.115"
$ws.Range("E39").Value = "  +20.66%  "
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("D41").Value = "17.93"
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "97.99"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").Value = "4.66"
$ws.Range("E46").Value = "  +15.61%  "
$ws.Range("E47").Value = "  -20.29%  "
$ws.Range("D48").Value = "2.49"
$ws.Range("E48").Value = "  +5.79%  "
$ws.Range("D49").Value = "1.294.41"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").Value = "2.92"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "6.93"
$ws.Range("E51").Value = "  -0.81%  "
